# Generate Report for Handback
#
# Updates the localization-status workbook to reflect that the 3e4e5c5f...
# and 4419680f... files have been handed back (in sync) for both the
# zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36a856b5d0221818b0f4b3b85f84895f47a917dc/e2e/3e4e5c5f-c122-49f3-a7bf-32de6ef3a284.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/36a856b5d0221818b0f4b3b85f84895f47a917dc/e2e/4419680f-1411-4177-8739-d2262446e549.md"

$name1 = "3e4e5c5f-c122-49f3-a7bf-32de6ef3a284.md"
$name2 = "4419680f-1411-4177-8739-d2262446e549.md"

# ---------------------------------------------------------------------
# Overview sheet: status text for both locales changes from
# "Ready for handoff" to "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.16666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

# Re-create the A2/A3 hyperlinks first so that the new I2/I3 hyperlinks
# are interleaved in the same order as the target workbook
# (A2, I2, A3, I3).
$wsZhCn.Range("A3").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlMd1, "", "", $name1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlMd1, "", "", $name1)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlMd2, "", "", $name2)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlMd2, "", "", $name2)

$wsZhCn.Range("J2").Value = "3e4e5c5f-c122-49f3-a7bf-32de6ef3a284.a7d851ce6afef47be524459f0d30b1dc47fd6890.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "4419680f-1411-4177-8739-d2262446e549.4d59612a54a84af0a380315a32ea46b63f8bc26f.zh-cn.xlf"

$wsZhCn.Range("K2").Value = "2016-08-24 02:29:15"
$wsZhCn.Range("K3").Value = "2016-08-24 02:29:15"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("A3").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlMd1, "", "", $name1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlMd1, "", "", $name1)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlMd2, "", "", $name2)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlMd2, "", "", $name2)

$wsDeDe.Range("J2").Value = "3e4e5c5f-c122-49f3-a7bf-32de6ef3a284.a7d851ce6afef47be524459f0d30b1dc47fd6890.de-de.xlf"
$wsDeDe.Range("J3").Value = "4419680f-1411-4177-8739-d2262446e549.4d59612a54a84af0a380315a32ea46b63f8bc26f.de-de.xlf"

$wsDeDe.Range("K2").Value = "2016-08-24 02:29:22"
$wsDeDe.Range("K3").Value = "2016-08-24 02:29:22"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.16666666666667
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16666666666667
